$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6397.5
$ws.Range("I32").Value = 5530
$ws.Range("K32").Value = 5530
$ws.Range("M32").Value = -5204

$ws.Range("H92").Value = 329.08694
$ws.Range("I92").Value = 250.875
$ws.Range("K92").Value = 250.875
$ws.Range("M92").Value = 997.125

$ws.Range("H103").Value = 20833858
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 22727810
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 68183430
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -68184602

$ws.Range("H117").Value = 66914
$ws.Range("I117").Value = 50000
$ws.Range("K117").Value = 50000
$ws.Range("M117").Value = -45411

$ws.Range("H132").Value = 3894.9333
$ws.Range("I132").Value = 1958.091
$ws.Range("J132").Value = 9221.25
$ws.Range("K132").Value = 5874.272999999999
$ws.Range("L132").Value = 27663.75
$ws.Range("M132").Value = -3344.272999999999
$ws.Range("N132").Value = -32723.75

$ws.Range("H134").Value = 181613.92
$ws.Range("I134").Value = 139738
$ws.Range("K134").Value = 139738
$ws.Range("M134").Value = -134668

$ws.Range("H138").Value = 2203.1667
$ws.Range("I138").Value = 1425.5264
$ws.Range("J138").Value = 5158.2
$ws.Range("K138").Value = 4276.5792
$ws.Range("L138").Value = 15474.6
$ws.Range("M138").Value = 863.4207999999999
$ws.Range("N138").Value = -25754.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2946
$ws.Range("I2").Value = 3091.25
$ws.Range("J2").Value = 2771.7
$ws.Range("K2").Value = 3091.25
$ws.Range("L2").Value = 2771.7
$ws.Range("M2").Value = -2978.25
$ws.Range("N2").Value = -2997.7

$ws.Range("H32").Value = 37946.19
$ws.Range("I32").Value = 21311.215
$ws.Range("J32").Value = 159143.86
$ws.Range("K32").Value = 21311.215
$ws.Range("L32").Value = 159143.86
$ws.Range("M32").Value = -21024.215
$ws.Range("N32").Value = -159717.86

$ws.Range("H61").Value = 1443.2084
$ws.Range("I61").Value = 991.75
$ws.Range("J61").Value = 2346.125
$ws.Range("K61").Value = 991.75
$ws.Range("L61").Value = 2346.125
$ws.Range("M61").Value = -779.75
$ws.Range("N61").Value = -2770.125

$ws.Range("H74").Value = 2317.8
$ws.Range("I74").Value = 1978
$ws.Range("J74").Value = 3252.25
$ws.Range("K74").Value = 1978
$ws.Range("L74").Value = 3252.25
$ws.Range("M74").Value = -1104
$ws.Range("N74").Value = -5000.25

$ws.Range("H77").Value = 2317.8
$ws.Range("I77").Value = 1978
$ws.Range("J77").Value = 3252.25
$ws.Range("K77").Value = 9890
$ws.Range("L77").Value = 16261.25
$ws.Range("M77").Value = -5522
$ws.Range("N77").Value = -24997.25

$ws.Range("H94").Value = 30330
$ws.Range("J94").Value = 30330
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -32132

$ws.Range("H97").Value = 10102727
$ws.Range("I97").Value = 12346841
$ws.Range("K97").Value = 12346841
$ws.Range("M97").Value = -12346345

$ws.Range("H116").Value = 2946
$ws.Range("I116").Value = 3091.25
$ws.Range("J116").Value = 2771.7
$ws.Range("K116").Value = 3091.25
$ws.Range("L116").Value = 2771.7
$ws.Range("M116").Value = -797.25
$ws.Range("N116").Value = -7359.7

$ws.Range("H131").Value = 79989
$ws.Range("J131").Value = 79989
$ws.Range("L131").Value = 79989
$ws.Range("N131").Value = -90069

$ws.Range("H136").Value = 1443.2084
$ws.Range("I136").Value = 991.75
$ws.Range("J136").Value = 2346.125
$ws.Range("K136").Value = 2975.25
$ws.Range("L136").Value = 7038.375
$ws.Range("M136").Value = -425.25
$ws.Range("N136").Value = -12138.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2946
$ws.Range("I3").Value = 3091.25
$ws.Range("J3").Value = 2771.7
$ws.Range("K3").Value = 3091.25
$ws.Range("L3").Value = 2771.7
$ws.Range("M3").Value = -2977.25
$ws.Range("N3").Value = -2999.7

$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 10000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -10470

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H86").Value = 2471.5557
$ws.Range("I86").Value = 1450.1666
$ws.Range("K86").Value = 1450.1666
$ws.Range("M86").Value = -327.1666

$ws.Range("H89").Value = 2471.5557
$ws.Range("I89").Value = 1450.1666
$ws.Range("K89").Value = 7250.833000000001
$ws.Range("M89").Value = -1634.833000000001

$ws.Range("H94").Value = 1310.4166
$ws.Range("I94").Value = 1039.2
$ws.Range("K94").Value = 1039.2
$ws.Range("M94").Value = -588.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1834.6666
$ws.Range("J16").Value = 1924
$ws.Range("L16").Value = 1924
$ws.Range("N16").Value = -2498

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H86").Value = 9015.5
$ws.Range("I86").Value = 8024.8
$ws.Range("K86").Value = 8024.8
$ws.Range("M86").Value = -6901.8

$ws.Range("H89").Value = 9015.5
$ws.Range("I89").Value = 8024.8
$ws.Range("K89").Value = 40124
$ws.Range("M89").Value = -34508

$ws.Range("H100").Value = 69250
$ws.Range("J100").Value = 69250
$ws.Range("L100").Value = 69250
$ws.Range("N100").Value = -71414

$ws.Range("H113").Value = 1834.6666
$ws.Range("J113").Value = 1924
$ws.Range("L113").Value = 1924
$ws.Range("N113").Value = -6264

$ws.Range("H130").Value = 74109
$ws.Range("J130").Value = 74109
$ws.Range("L130").Value = 74109
$ws.Range("N130").Value = -84149

$ws.Range("H131").Value = 49999
$ws.Range("J131").Value = 49999
$ws.Range("L131").Value = 49999
$ws.Range("N131").Value = -60079

$ws.Range("H132").Value = 2136.6667
$ws.Range("I132").Value = 1884.3334
$ws.Range("J132").Value = 4660
$ws.Range("K132").Value = 5653.0002
$ws.Range("L132").Value = 13980
$ws.Range("M132").Value = -3123.0002
$ws.Range("N132").Value = -19040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 450
$ws.Range("I2").Value = 450
$ws.Range("K2").Value = 2700
$ws.Range("M2").Value = -2587

$ws.Range("H9").Value = 12707.286
$ws.Range("J9").Value = 17000.4
$ws.Range("L9").Value = 51001.2
$ws.Range("N9").Value = -51449.2

$ws.Range("H108").Value = 2931.75
$ws.Range("I108").Value = 1779.2858
$ws.Range("K108").Value = 5337.857400000001
$ws.Range("M108").Value = -2457.857400000001

$ws.Range("H114").Value = 11765518
$ws.Range("I114").Value = 22222992
$ws.Range("J114").Value = 861.125
$ws.Range("K114").Value = 66668976
$ws.Range("L114").Value = 2583.375
$ws.Range("M114").Value = -66665722
$ws.Range("N114").Value = -9091.375

$ws.Range("H116").Value = 9149.786
$ws.Range("I116").Value = 9084.538
$ws.Range("K116").Value = 27253.614
$ws.Range("M116").Value = -23811.614

$ws.Range("H121").Value = 17598056
$ws.Range("J121").Value = 2192.5
$ws.Range("L121").Value = 6577.5
$ws.Range("N121").Value = -9197.5

$ws.Range("H129").Value = 112783.945
$ws.Range("J129").Value = 3433
$ws.Range("L129").Value = 10299
$ws.Range("N129").Value = -20299

$ws.Range("H131").Value = 53553.19
$ws.Range("I131").Value = 100538
$ws.Range("K131").Value = 301614
$ws.Range("M131").Value = -296574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H97").Value = 22830.533
$ws.Range("I97").Value = 36985.445
$ws.Range("K97").Value = 36985.445
$ws.Range("M97").Value = -36489.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 40000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 40000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 40000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -41322

$ws.Range("H132").Value = 558476.7
$ws.Range("I132").Value = 683858
$ws.Range("J132").Value = 6799
$ws.Range("K132").Value = 2051574
$ws.Range("L132").Value = 20397
$ws.Range("M132").Value = -2049044
$ws.Range("N132").Value = -25457

$ws.Range("H135").Value = 16803.572
$ws.Range("I135").Value = 12326.923
$ws.Range("K135").Value = 12326.923
$ws.Range("M135").Value = -7256.923000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178

$ws.Range("H124").Value = 61474.332
$ws.Range("J124").Value = 61474.332
$ws.Range("L124").Value = 61474.332
$ws.Range("N124").Value = -71294.33199999999

$ws.Range("H126").Value = 2470.5264
$ws.Range("I126").Value = 1964.7693
$ws.Range("K126").Value = 5894.3079
$ws.Range("M126").Value = -3424.3079

